# Generate Report for Archive
# Four files (rows 2-5 in each sheet) move from "Ready for handoff" to
# "In Translation": 0de4acb4, 4b8c1883, 712ccd48, c0dbaea3.
# Rows 6-8 (e9df6f50, f5df241f, f8979864) stay "Ready for handoff".

$wb = $excel.ActiveWorkbook

# --- zh-cn / de-de sheets: Status column is column C, data rows 2-5 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 5; $r++) {
        $ws.Cells.Item($r, 3).Value = "In Translation"
    }
}

# --- Overview sheet: zh-cn column is E, de-de column is F, data rows 2-5 ---
$overview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 5; $r++) {
    $overview.Cells.Item($r, 5).Value = "In Translation"
    $overview.Cells.Item($r, 6).Value = "In Translation"
}
